{"js": "// Change the East Asian / Complex Script fallback fonts used by the\n// document's paragraph styles from \"DejaVu Sans\" to \"Tahoma\" (East Asian)\n// and make sure the Complex Script font is explicitly pinned to\n// \"DejaVu Sans\" on the styles that did not previously carry an explicit\n// <w:rFonts> override.\n//\n// Word.Style.font.nameFarEast      -> <w:rFonts w:eastAsia=\"\u2026\"/>\n// Word.Style.font.nameBidirectional -> <w:rFonts w:cs=\"\u2026\"/>\n\nconst styles = context.document.getStyles();\nstyles.load(\"items/nameLocal\");\nawait context.sync();\n\n// Build a lookup by style name so we don't depend on collection order.\nconst byName = {};\nfor (const s of styles.items) {\n  byName[s.nameLocal] = s;\n}\n\n// \"Normal\" and \"Heading\" already define an explicit eastAsia font\n// (DejaVu Sans) -> switch it to Tahoma.\nif (byName[\"Normal\"]) {\n  byName[\"Normal\"].font.nameFarEast = \"Tahoma\";\n}\nif (byName[\"Heading\"]) {\n  byName[\"Heading\"].font.nameFarEast = \"Tahoma\";\n}\n\n// \"List\", \"Caption\" and \"Index\" had no direct rFonts override (they\n// inherited everything) -> pin their complex-script font explicitly to\n// \"DejaVu Sans\".\nif (byName[\"List\"]) {\n  byName[\"List\"].font.nameBidirectional = \"DejaVu Sans\";\n}\nif (byName[\"Caption\"]) {\n  byName[\"Caption\"].font.nameBidirectional = \"DejaVu Sans\";\n}\nif (byName[\"Index\"]) {\n  byName[\"Index\"].font.nameBidirectional = \"DejaVu Sans\";\n}\n\nawait context.sync();\n", "ps1": "# Change the East Asian / Complex Script fallback fonts used by the\n# document's paragraph styles from \"DejaVu Sans\" to \"Tahoma\" (East Asian)\n# and make sure the Complex Script font is explicitly pinned to\n# \"DejaVu Sans\" on the styles that did not previously carry an explicit\n# <w:rFonts> override.\n#\n# Style.Font.NameFarEast -> <w:rFonts w:eastAsia=\"\u2026\"/>\n# Style.Font.NameBi      -> <w:rFonts w:cs=\"\u2026\"/>\n\n$d = $word.ActiveDocument\n\n# \"Normal\" and \"Heading\" already define an explicit eastAsia font\n# (DejaVu Sans) -> switch it to Tahoma.\n$d.Styles.Item(\"Normal\").Font.NameFarEast = \"Tahoma\"\n$d.Styles.Item(\"Heading\").Font.NameFarEast = \"Tahoma\"\n\n# \"List\", \"Caption\" and \"Index\" had no direct rFonts override (they\n# inherited everything) -> pin their complex-script font explicitly to\n# \"DejaVu Sans\".\n$d.Styles.Item(\"List\").Font.NameBi = \"DejaVu Sans\"\n$d.Styles.Item(\"Caption\").Font.NameBi = \"DejaVu Sans\"\n$d.Styles.Item(\"Index\").Font.NameBi = \"DejaVu Sans\"\n"}
